$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stage Price (D) column values as text via helper column Z, then paste-special values ---
$ws.Range("Z2").Value = "'27.944.86"
$ws.Range("Z3").Value = "'1.866.52"
$ws.Range("Z4").Value = "'1.000"
$ws.Range("Z5").Value = "'317.54"
$ws.Range("Z6").Value = "'0.9997"
$ws.Range("Z7").Value = "'0.4369"
$ws.Range("Z8").Value = "'0.3722"
$ws.Range("Z9").Value = "'0.07490"
$ws.Range("Z10").Value = "'0.9370"
$ws.Range("Z11").Value = "'21.37"
$ws.Range("Z12").Value = "'1.875.10"
$ws.Range("Z13").Value = "'6.744"
$ws.Range("Z14").Value = "'5.451"
$ws.Range("Z15").Value = "'0.06852"
$ws.Range("Z16").Value = "'1.002"
$ws.Range("Z17").Value = "'81.61"
$ws.Range("Z18").Value = "'0.000009066"
$ws.Range("Z19").Value = "'0.9988"
$ws.Range("Z20").Value = "'15.92"
$ws.Range("Z21").Value = "'27.932.56"
$ws.Range("Z22").Value = "'5.134"
$ws.Range("Z23").Value = "'11.04"
$ws.Range("Z24").Value = "'2.097.20"
$ws.Range("Z25").Value = "'2.004"
$ws.Range("Z26").Value = "'154.28"
$ws.Range("Z27").Value = "'18.48"
$ws.Range("Z28").Value = "'5.498"
$ws.Range("Z29").Value = "'113.22"
$ws.Range("Z30").Value = "'1.721"
$ws.Range("Z31").Value = "'0.09032"
$ws.Range("Z32").Value = "'0.8241"
$ws.Range("Z33").Value = "'4.824"
$ws.Range("Z34").Value = "'1.178"
$ws.Range("Z35").Value = "'2.978"
$ws.Range("Z36").Value = "'0.9999"
$ws.Range("Z37").Value = "'0.05510"
$ws.Range("Z38").Value = "'1.120"
$ws.Range("Z39").Value = "'0.01982"
$ws.Range("Z40").Value = "'2.955"
$ws.Range("Z41").Value = "'0.5269"
$ws.Range("Z42").Value = "'7.053"
$ws.Range("Z43").Value = "'0.1706"
$ws.Range("Z44").Value = "'8.806"
$ws.Range("Z45").Value = "'0.06771"
$ws.Range("Z46").Value = "'0.4909"
$ws.Range("Z47").Value = "'10.71"
$ws.Range("Z48").Value = "'107.31"
$ws.Range("Z49").Value = "'1.682"
$ws.Range("Z50").Value = "'0.9992"
$ws.Range("Z51").Value = "'1.889"
$ws.Range("Z2:Z51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$ws.Range("Z2:Z51").Clear()

# --- Volume(1h) (E) column: plain text assignment (always contains "%" so never numeric) ---
$ws.Range("E2").Value = "  -3.58%  "
$ws.Range("E3").Value = "  -2.62%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -4.98%  "
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("E10").Value = "  -4.48%  "
$ws.Range("E11").Value = "  -3.73%  "
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("E14").Value = "  -4.17%  "
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("E18").Value = "  -4.45%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  -4.51%  "
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("E22").Value = "  -3.65%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  -3.66%  "
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("E30").Value = "  -7.84%  "
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("E32").Value = "  -4.83%  "
$ws.Range("E33").Value = "  -5.66%  "
$ws.Range("E34").Value = "  -5.79%  "
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("E41").Value = "  -4.24%  "
$ws.Range("E42").Value = "  -6.16%  "
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("E44").Value = "  -6.23%  "
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("E46").Value = "  -5.47%  "
$ws.Range("E47").Value = "  -4.53%  "
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("E49").Value = "  -5.62%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").Value = "  -13.34%  "

# --- Row 50/51 Coin name + Link swap ---
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
